# ECP-962: updates turnover import
# Renames the TurnoverImport header row to Title-Case business names,
# adds four new trailing columns (Lease Name, Gross/Net Amount Previous Year,
# Purchase Count Previous Year), restyles the header row with a new font,
# bumps the header row height, re-sizes several columns, and leaves a
# formatted-but-empty marker cell at C6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Rename the existing headers (A1:N1) to their Title Case equivalents
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Lease Reference"
$ws.Range("B1").Value = "Unit Reference"
$ws.Range("C1").Value = "Occupancy Start Date"
$ws.Range("D1").Value = "Date"
$ws.Range("E1").Value = "Type"
$ws.Range("F1").Value = "Frequency"
$ws.Range("G1").Value = "Gross Amount"
$ws.Range("H1").Value = "Net Amount"
$ws.Range("I1").Value = "Non Comparable Flag"
$ws.Range("J1").Value = "Purchase Count"
$ws.Range("K1").Value = "Comments"
$ws.Range("L1").Value = "Currency"
$ws.Range("M1").Value = "Reported By"
$ws.Range("N1").Value = "Reported At"

# ---------------------------------------------------------------------
# 2. Add the four new trailing headers (O1:R1)
# ---------------------------------------------------------------------
$ws.Range("O1").Value = "Lease Name"
$ws.Range("P1").Value = "Gross Amount Previous Year"
$ws.Range("Q1").Value = "Net Amount Previous Year"
$ws.Range("R1").Value = "Purchase Count Previous Year"

# ---------------------------------------------------------------------
# 3. Re-style the whole header row with the new font (Calibri 11) and
#    a taller row height, matching the "previous year" marker cell C6.
# ---------------------------------------------------------------------
$headerRow = $ws.Range("A1:R1")
$headerRow.Font.Size = 11
$headerRow.Font.Name = "Calibri"
$ws.Rows.Item(1).RowHeight = 14

$ws.Range("C6").Font.Size = 11
$ws.Range("C6").Font.Name = "Calibri"

# ---------------------------------------------------------------------
# 4. Resize columns to fit the new header text
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 11.498697916666666   # B  12.33203125
$ws.Columns.Item(3).ColumnWidth = 15.498697916666666   # C  16.33203125
$ws.Columns.Item(6).ColumnWidth = 9.498697916666666    # F  10.33203125
$ws.Columns.Item(7).ColumnWidth = 10.998697916666666   # G  11.83203125
$ws.Columns.Item(8).ColumnWidth = 9.998697916666666    # H  10.83203125
$ws.Columns.Item(9).ColumnWidth = 15.666666666666666   # I  16.5
$ws.Columns.Item(10).ColumnWidth = 19.998697916666668  # J  20.83203125
$ws.Columns.Item(11).ColumnWidth = 36.666666666666664  # K  37.5
$ws.Columns.Item(12).ColumnWidth = 16.498697916666668  # L  17.33203125
$ws.Columns.Item(13).ColumnWidth = 10.166666666666666  # M  11
$ws.Columns.Item(14).ColumnWidth = 18.330729166666668  # N  19.1640625
$ws.Columns.Item(15).ColumnWidth = 9.498697916666666   # O  10.33203125
$ws.Columns.Item(16).ColumnWidth = 21.498697916666668  # P  22.33203125
$ws.Columns.Item(17).ColumnWidth = 19.998697916666668  # Q  20.83203125
$ws.Columns.Item(18).ColumnWidth = 23.330729166666668  # R  24.1640625

Write-Output "TurnoverImport headers updated"
